$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(4, 2).Value = 8530411
$ws.Cells.Item(4, 3).Value = 9461
$ws.Cells.Item(4, 4).Value = 5551127
$ws.Cells.Item(4, 5).Value = 2752866
$ws.Cells.Item(4, 7).Value = 234
$ws.Cells.Item(4, 8).Value = 226418
$ws.Cells.Item(5, 2).Value = 7701365
$ws.Cells.Item(5, 3).Value = 52207
$ws.Cells.Item(5, 4).Value = 6867988
$ws.Cells.Item(5, 5).Value = 716792
$ws.Cells.Item(5, 7).Value = 635
$ws.Cells.Item(5, 8).Value = 116585
$ws.Cells.Item(6, 2).Value = 5276942
$ws.Cells.Item(6, 3).Value = 2125
$ws.Cells.Item(6, 5).Value = 400443
$ws.Cells.Item(6, 7).Value = 18
$ws.Cells.Item(6, 8).Value = 154906
$ws.Cells.Item(14, 2).Value = 789229
$ws.Cells.Item(14, 3).Value = 26688
$ws.Cells.Item(14, 7).Value = 191
$ws.Cells.Item(14, 8).Value = 44158
$ws.Cells.Item(21, 2).Value = 385342
$ws.Cells.Item(21, 3).Value = 4444
$ws.Cells.Item(21, 5).Value = 73276
$ws.Cells.Item(21, 7).Value = 11
$ws.Cells.Item(21, 8).Value = 9966
$ws.Cells.Item(24, 2).Value = 353426
$ws.Cells.Item(24, 3).Value = 2013
$ws.Cells.Item(24, 4).Value = 308446
$ws.Cells.Item(24, 5).Value = 35467
$ws.Cells.Item(24, 7).Value = 68
$ws.Cells.Item(24, 8).Value = 9513
$ws.Cells.Item(28, 2).Value = 306649
$ws.Cells.Item(28, 3).Value = 487
$ws.Cells.Item(28, 4).Value = 283207
$ws.Cells.Item(28, 5).Value = 21164
$ws.Cells.Item(31, 2).Value = 205557
$ws.Cells.Item(31, 3).Value = 1869
$ws.Cells.Item(31, 4).Value = 173269
$ws.Cells.Item(31, 5).Value = 22465
$ws.Cells.Item(31, 7).Value = 29
$ws.Cells.Item(31, 8).Value = 9823
$ws.Cells.Item(32, 1).Value = "Chequia"
$ws.Cells.Item(32, 2).Value = 202787
$ws.Cells.Item(32, 3).Value = 8841
$ws.Cells.Item(32, 4).Value = 83097
$ws.Cells.Item(32, 5).Value = 117971
$ws.Cells.Item(32, 7).Value = 100
$ws.Cells.Item(32, 8).Value = 1719
$ws.Cells.Item(33, 1).Value = "Polonia"
$ws.Cells.Item(33, 2).Value = 202579
$ws.Cells.Item(33, 3).Value = 10040
$ws.Cells.Item(33, 4).Value = 98884
$ws.Cells.Item(33, 5).Value = 99844
$ws.Cells.Item(33, 7).Value = 130
$ws.Cells.Item(33, 8).Value = 3851
$ws.Cells.Item(53, 2).Value = 91118
$ws.Cells.Item(53, 3).Value = 628
$ws.Cells.Item(53, 4).Value = 44506
$ws.Cells.Item(53, 5).Value = 45228
$ws.Cells.Item(53, 7).Value = 13
$ws.Cells.Item(53, 8).Value = 1384
$ws.Cells.Item(67, 2).Value = 55081
$ws.Cells.Item(67, 3).Value = 252
$ws.Cells.Item(67, 4).Value = 38482
$ws.Cells.Item(67, 5).Value = 14719
$ws.Cells.Item(67, 7).Value = 7
$ws.Cells.Item(67, 8).Value = 1880
$ws.Cells.Item(69, 2).Value = 53422
$ws.Cells.Item(69, 3).Value = 1166
$ws.Cells.Item(69, 5).Value = 28190
$ws.Cells.Item(69, 7).Value = 3
$ws.Cells.Item(69, 8).Value = 1868
$ws.Cells.Item(77, 2).Value = 43620
$ws.Cells.Item(77, 3).Value = 2648
$ws.Cells.Item(77, 4).Value = 7223
$ws.Cells.Item(77, 5).Value = 35954
$ws.Cells.Item(77, 7).Value = 29
$ws.Cells.Item(77, 8).Value = 443
$ws.Cells.Item(88, 2).Value = 27334
$ws.Cells.Item(88, 3).Value = 865
$ws.Cells.Item(88, 5).Value = 16811
$ws.Cells.Item(88, 7).Value = 6
$ws.Cells.Item(88, 8).Value = 534
$ws.Cells.Item(90, 2).Value = 24836
$ws.Cells.Item(90, 3).Value = 640
$ws.Cells.Item(90, 4).Value = 17905
$ws.Cells.Item(90, 5).Value = 6069
$ws.Cells.Item(90, 7).Value = 12
$ws.Cells.Item(90, 8).Value = 862
$ws.Cells.Item(98, 1).Value = "Montenegro"
$ws.Cells.Item(98, 2).Value = 16069
$ws.Cells.Item(98, 3).Value = 177
$ws.Cells.Item(98, 4).Value = 11815
$ws.Cells.Item(98, 5).Value = 4004
$ws.Cells.Item(98, 7).Value = 3
$ws.Cells.Item(98, 8).Value = 250
$ws.Cells.Item(99, 1).Value = "Zambia"
$ws.Cells.Item(99, 2).Value = 16000
$ws.Cells.Item(99, 3).Value = 18
$ws.Cells.Item(99, 4).Value = 15168
$ws.Cells.Item(99, 5).Value = 486
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 346
$ws.Cells.Item(100, 1).Value = "Eslovenia"
$ws.Cells.Item(100, 2).Value = 15982
$ws.Cells.Item(100, 3).Value = 1503
$ws.Cells.Item(100, 4).Value = 6922
$ws.Cells.Item(100, 5).Value = 8860
$ws.Cells.Item(100, 7).Value = 8
$ws.Cells.Item(100, 8).Value = 200
$ws.Cells.Item(105, 1).Value = "Luxemburgo"
$ws.Cells.Item(105, 2).Value = 11671
$ws.Cells.Item(105, 3).Value = 430
$ws.Cells.Item(105, 4).Value = 8473
$ws.Cells.Item(105, 5).Value = 3060
$ws.Cells.Item(105, 7).Value = 2
$ws.Cells.Item(105, 8).Value = 138
$ws.Cells.Item(106, 1).Value = "Guinea"
$ws.Cells.Item(106, 2).Value = 11599
$ws.Cells.Item(106, 3).Value = 61
$ws.Cells.Item(106, 4).Value = 10461
$ws.Cells.Item(106, 5).Value = 1068
$ws.Cells.Item(106, 8).Value = 70
$ws.Cells.Item(107, 1).Value = "Maldivas"
$ws.Cells.Item(107, 2).Value = 11271
$ws.Cells.Item(107, 4).Value = 10234
$ws.Cells.Item(107, 5).Value = 1000
$ws.Cells.Item(107, 8).Value = 37
$ws.Cells.Item(124, 2).Value = 5977
$ws.Cells.Item(124, 3).Value = 166
$ws.Cells.Item(124, 5).Value = 2463
$ws.Cells.Item(137, 1).Value = "Reunion"
$ws.Cells.Item(137, 2).Value = 5015
$ws.Cells.Item(137, 3).Value = 94
$ws.Cells.Item(137, 4).Value = 4445
$ws.Cells.Item(137, 5).Value = 551
$ws.Cells.Item(137, 7).Value = 2
$ws.Cells.Item(137, 8).Value = 19
$ws.Cells.Item(138, 1).Value = "Ruanda"
$ws.Cells.Item(138, 2).Value = 4996
$ws.Cells.Item(138, 4).Value = 4797
$ws.Cells.Item(138, 5).Value = 165
$ws.Cells.Item(138, 8).Value = 34
$ws.Cells.Item(143, 1).Value = "Mayotte"
$ws.Cells.Item(143, 2).Value = 4203
$ws.Cells.Item(143, 4).Value = 2964
$ws.Cells.Item(143, 5).Value = 1195
$ws.Cells.Item(143, 7).Value = 1
$ws.Cells.Item(143, 8).Value = 44
$ws.Cells.Item(144, 1).Value = "Estonia"
$ws.Cells.Item(144, 2).Value = 4171
$ws.Cells.Item(144, 3).Value = 44
$ws.Cells.Item(144, 4).Value = 3334
$ws.Cells.Item(144, 5).Value = 766
$ws.Cells.Item(144, 8).Value = 71
$ws.Cells.Item(146, 1).Value = "Principado de Andorra"
$ws.Cells.Item(146, 2).Value = 3811
$ws.Cells.Item(146, 4).Value = 2470
$ws.Cells.Item(146, 5).Value = 1278
$ws.Cells.Item(146, 7).Value = 1
$ws.Cells.Item(146, 8).Value = 63
$ws.Cells.Item(147, 1).Value = "Letonia"
$ws.Cells.Item(147, 2).Value = 3797
$ws.Cells.Item(147, 3).Value = 188
$ws.Cells.Item(147, 4).Value = 1341
$ws.Cells.Item(147, 5).Value = 2409
$ws.Cells.Item(147, 8).Value = 47
$ws.Cells.Item(148, 1).Value = "Guyana"
$ws.Cells.Item(148, 2).Value = 3796
$ws.Cells.Item(148, 3).Value = 0
$ws.Cells.Item(148, 4).Value = 2796
$ws.Cells.Item(148, 5).Value = 886
$ws.Cells.Item(148, 8).Value = 114
$ws.Cells.Item(149, 1).Value = "Tailandia"
$ws.Cells.Item(149, 2).Value = 3709
$ws.Cells.Item(149, 3).Value = 9
$ws.Cells.Item(149, 4).Value = 3495
$ws.Cells.Item(149, 5).Value = 155
$ws.Cells.Item(149, 8).Value = 59
$ws.Cells.Item(150, 1).Value = "Gambia"
$ws.Cells.Item(150, 2).Value = 3655
$ws.Cells.Item(150, 4).Value = 2658
$ws.Cells.Item(150, 5).Value = 879
$ws.Cells.Item(150, 8).Value = 118
$ws.Cells.Item(152, 1).Value = "Republica de Chipre"
$ws.Cells.Item(152, 2).Value = 2966
$ws.Cells.Item(152, 3).Value = 127
$ws.Cells.Item(152, 4).Value = 1444
$ws.Cells.Item(152, 5).Value = 1497
$ws.Cells.Item(152, 8).Value = 25
$ws.Cells.Item(153, 1).Value = "Belice"
$ws.Cells.Item(153, 2).Value = 2886
$ws.Cells.Item(153, 3).Value = 53
$ws.Cells.Item(153, 4).Value = 1727
$ws.Cells.Item(153, 5).Value = 1114
$ws.Cells.Item(153, 8).Value = 45
$ws.Cells.Item(154, 1).Value = "Sudan del Sur"
$ws.Cells.Item(154, 2).Value = 2847
$ws.Cells.Item(154, 4).Value = 1290
$ws.Cells.Item(154, 5).Value = 1502
$ws.Cells.Item(154, 8).Value = 55
$ws.Cells.Item(159, 2).Value = 2337
$ws.Cells.Item(159, 3).Value = 1
$ws.Cells.Item(159, 4).Value = 1771
$ws.Cells.Item(159, 5).Value = 493
$ws.Cells.Item(192, 4).Value = 205
$ws.Cells.Item(192, 5).Value = 10
$ws.Cells.Item(203, 2).Value = 38
$ws.Cells.Item(203, 3).Value = 2
$ws.Cells.Item(203, 5).Value = 11
